$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "number of convolution layers" value in C2
$ws.Range("C2").Value = "[ [32,64,128,160,180,256, ], [64,80,128,256,270,364,  ], [80,100, 128,256,512,712, ], [128, 170,256,512,1024,2048,],  ]"

# Update column widths for B and C
$ws.Columns.Item(2).ColumnWidth = 84.7109375
$ws.Columns.Item(3).ColumnWidth = 102.28125

# Update the view's top-left visible cell from C1 to B1
$ws.Application.ActiveWindow.ScrollColumn = 2
